$d = $word.ActiveDocument

# 1. Créditos-aula: 4 -> 2
$d.Content.Find.Execute("Créditos-aula: 4", $true, $false, $false, $false, $false, $true, 1, $false, "Créditos-aula: 2", 2)

# 2. Carga horária: 60 h -> 30 h
$d.Content.Find.Execute("Carga horária: 60 h", $true, $false, $false, $false, $false, $true, 1, $false, "Carga horária: 30 h", 2)

# 3. Ativação: 01/01/2012 -> 01/01/2023
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2)

# 4. Curso (semestre ideal): EF (7) -> EF (5)
$d.Content.Find.Execute("Curso (semestre ideal): EF (7)", $true, $false, $false, $false, $false, $true, 1, $false, "Curso (semestre ideal): EF (5)", 2)

# 5. Replace objectives paragraph (PT) text
$oldObjPt = "Apresentar noções de trocas de calor, mediante estudo dos mecanismos básicos. Capacitar o aluno a modelar e resolver problemas de interesse em transferência de calor, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução."
$newObjPt = "Introdução de conceitos relacionados com taxa e fluxo de quantidade de movimento, calor e massa. Capacitar o aluno a modelar e resolver problemas de interesse em fenômenos de transporte, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução."
$d.Content.Find.Execute($oldObjPt, $true, $false, $false, $false, $false, $true, 1, $false, $newObjPt, 2)

# 6. Insert new italic English paragraph after the objectives paragraph
$objPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Introdução de conceitos relacionados com taxa e fluxo*") {
        $objPara = $p
        break
    }
}
$objPara.Range.InsertParagraphAfter()
$found = $false
$objEnPara = $null
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $objEnPara = $p
        break
    }
    if ($p.Range.Text -like "*Introdução de conceitos relacionados com taxa e fluxo*") {
        $found = $true
    }
}
$objEnRange = $objEnPara.Range
$objEnRange.MoveEnd(1, -1)
$objEnRange.Text = "Introduction of concepts related to rate and flow of momentum, heat and mass. Enable the student to model and solve problems of interest in transport phenomena, with appropriate choice of hypotheses and application of corresponding solution tools."
$objEnRange.Font.Italic = 1

# 7. Remove "519033 - Carlos Yujiro Shigue" run (with its line break)
$d.Content.Find.Execute("519033 - Carlos Yujiro Shigue" + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 8. Replace "Programa resumido" paragraph (PT) text
$oldResumoPt = "Introdução à transferência de calor. Condução de calor em regime permanente e em regime transiente. Convecção forçada em dutos e sobre corpos; convecção natural. Transferência de calor por radiação térmica. Transferência de calor com mudança de fase. Transferência de massa."
$newResumoPt = "Introdução à transferência de calor. Condução de calor em regime permanente e em regime transiente. Transferência de calor por convecção livre e forçada. Transferência de calor por radiação térmica. Transferência de calor com mudança de fase. Transferência de massa. Exemplos de aplicação."
$d.Content.Find.Execute($oldResumoPt, $true, $false, $false, $false, $false, $true, 1, $false, $newResumoPt, 2)

# 9. Insert new italic English paragraph after "Programa resumido" paragraph
$resumoPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Transferência de calor por convecção livre e forçada*") {
        $resumoPara = $p
        break
    }
}
$resumoPara.Range.InsertParagraphAfter()
$found = $false
$resumoEnPara = $null
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $resumoEnPara = $p
        break
    }
    if ($p.Range.Text -like "*Transferência de calor por convecção livre e forçada*") {
        $found = $true
    }
}
$resumoEnRange = $resumoEnPara.Range
$resumoEnRange.MoveEnd(1, -1)
$resumoEnRange.Text = "Introduction to heat transfer. Heat conduction in steady state and in transient regime. Free and forced convection heat transfer. Heat transfer by thermal radiation. Heat transfer with phase change. Mass transference. Application examples."
$resumoEnRange.Font.Italic = 1

# 10. Replace "Programa" paragraph: merge bulleted lines into one run with "•" separators
$programaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Transferência de calor por condução: transferência de calor unidimensional*") {
        $programaPara = $p
        break
    }
}
$newProgramaPt = "• Transferência de calor por condução: transferência de calor unidimensional em regime permanente. Equação de Fourier. Condutividade térmica. • Transferência de calor unidimensional em regime permanente com contornos convectivos. Lei de Newton do resfriamento. • Condução de calor em regime transiente. Difusividade térmica. Número de Biot. • Analogia entre transferência de calor e circuitos elétricos: conceitos de resistência e capacitância térmicas. • Transferência de calor 2D e 3D em regime transiente. • Transferência de calor por convecção. Convecção livre. Parâmetros de similiaridade. Número de Rayleigh. Convecção forçada. Teoria da camada limite. Número de Prandtl e número de Nusselt. • Transferência de calor por radiação. Radiação do corpo negro. Propriedades da radiação. Fator de forma da radiação. • Transferência de calor com mudança de fase: ebulição e condensação."
$programaRange = $programaPara.Range
$programaRange.MoveEnd(1, -1)
$programaRange.Text = $newProgramaPt

# 11. Insert new italic English paragraph after "Programa" paragraph
$programaPtPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*condutividade térmica*" -or $p.Range.Text -like "*Condutividade térmica*") {
        $programaPtPara = $p
        break
    }
}
$programaPtPara.Range.InsertParagraphAfter()
$found = $false
$programaEnPara = $null
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $programaEnPara = $p
        break
    }
    if ($p.Range.Text -like "*Condutividade térmica*") {
        $found = $true
    }
}
$programaEnRange = $programaEnPara.Range
$programaEnRange.MoveEnd(1, -1)
$programaEnRange.Text = "• Heat transfer by conduction: one-dimensional heat transfer in steady state. Fourier equation. Thermal conductivity. • One-dimensional heat transfer in steady state with convective contours. Newton's Law of Cooling. • Transient heat conduction. Thermal diffusivity. Biot number. • Analogy between heat transfer and electrical circuits: concepts of thermal resistance and capacitance. • Transient 2D and 3D heat transfer. • Convection heat transfer. Free convection. Similarity parameters. Rayleigh number. Forced convection. Boundary layer theory. Prandtl number and Nusselt number. • Radiation heat transfer. Blackbody radiation. Radiation properties. Radiation form factor. • Phase change heat transfer: boiling and condensation."
$programaEnRange.Font.Italic = 1
